$wb = $excel.ActiveWorkbook

# Sheet R1
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3929:40:22"
$ws.Range("G3").Value = "69:13:00"

# Sheet R2
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12111:04:03"
$ws.Range("G3").Value = "3240:47:32"
$ws.Range("G4").Value = "478:59:06"

# Sheet R4
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2956:53:52"
$ws.Range("G3").Value = "184:06:07"
$ws.Range("G4").Value = "72:18:32"
$ws.Range("G5").Value = "69:56:05"

# Sheet R5
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "430:52:51"

# Sheet R6
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "71:25:09"
